# Auto-generated edit script: adds 2022-Q4 sheet + updates 总计 summary sheet
$wb = $excel.ActiveWorkbook

# ---- 1. Insert a new worksheet named '2022-Q4' right before '2022-Q3' ----
$refSheet = $wb.Worksheets.Item('2022-Q3')
$q4 = $wb.Worksheets.Add($refSheet)
$q4.Name = '2022-Q4'

# Header row (bold/boxed style copied from the equivalent cell on the
# neighbouring quarter sheet so the new sheet matches the existing look).
$refSheet.Cells.Item(1,2).Copy($q4.Cells.Item(1,2))
$q4.Cells.Item(1,2).Value = '基金代码'
$refSheet.Cells.Item(1,3).Copy($q4.Cells.Item(1,3))
$q4.Cells.Item(1,3).Value = '基金名称'
$refSheet.Cells.Item(1,4).Copy($q4.Cells.Item(1,4))
$q4.Cells.Item(1,4).Value = '基金规模'
$refSheet.Cells.Item(1,5).Copy($q4.Cells.Item(1,5))
$q4.Cells.Item(1,5).Value = '股票总仓位'
$refSheet.Cells.Item(1,6).Copy($q4.Cells.Item(1,6))
$q4.Cells.Item(1,6).Value = '仓位占比'
$refSheet.Cells.Item(1,7).Copy($q4.Cells.Item(1,7))
$q4.Cells.Item(1,7).Value = '持有市值(亿元)'
$refSheet.Cells.Item(1,8).Copy($q4.Cells.Item(1,8))
$q4.Cells.Item(1,8).Value = '仓位排名'

# Data rows 2..25
# row 2: 000979
$refSheet.Cells.Item(2,1).Copy($q4.Cells.Item(2,1))
$q4.Cells.Item(2,1).Value = 0
$q4.Cells.Item(2,2).Value = '''000979'
$q4.Cells.Item(2,3).Value = '景顺长城沪港深精选股票'
$q4.Cells.Item(2,4).Value = '''25.01'
$q4.Cells.Item(2,5).Value = '''92.37'
$q4.Cells.Item(2,6).Value = '''8.02'
$q4.Cells.Item(2,7).Value = '''2.0058'
$q4.Cells.Item(2,8).Value = 5

# row 3: 260112
$refSheet.Cells.Item(3,1).Copy($q4.Cells.Item(3,1))
$q4.Cells.Item(3,1).Value = 1
$q4.Cells.Item(3,2).Value = '''260112'
$q4.Cells.Item(3,3).Value = '景顺长城能源基建混合A'
$q4.Cells.Item(3,4).Value = '''22.88'
$q4.Cells.Item(3,5).Value = '''85.17'
$q4.Cells.Item(3,6).Value = '''8.76'
$q4.Cells.Item(3,7).Value = '''2.0043'
$q4.Cells.Item(3,8).Value = 2

# row 4: 008850
$refSheet.Cells.Item(4,1).Copy($q4.Cells.Item(4,1))
$q4.Cells.Item(4,1).Value = 2
$q4.Cells.Item(4,2).Value = '''008850'
$q4.Cells.Item(4,3).Value = '景顺长城价值稳进三年定期开放灵活配置混合'
$q4.Cells.Item(4,4).Value = '''18.18'
$q4.Cells.Item(4,5).Value = '''97.45'
$q4.Cells.Item(4,6).Value = '''7.89'
$q4.Cells.Item(4,7).Value = '''1.4344'
$q4.Cells.Item(4,8).Value = 4

# row 5: 008715
$refSheet.Cells.Item(5,1).Copy($q4.Cells.Item(5,1))
$q4.Cells.Item(5,1).Value = 3
$q4.Cells.Item(5,2).Value = '''008715'
$q4.Cells.Item(5,3).Value = '景顺长城价值驱动一年持有期灵活配置混合'
$q4.Cells.Item(5,4).Value = '''6.75'
$q4.Cells.Item(5,5).Value = '''92.52'
$q4.Cells.Item(5,6).Value = '''9.94'
$q4.Cells.Item(5,7).Value = '''0.6710'
$q4.Cells.Item(5,8).Value = 1

# row 6: 009098
$refSheet.Cells.Item(6,1).Copy($q4.Cells.Item(6,1))
$q4.Cells.Item(6,1).Value = 4
$q4.Cells.Item(6,2).Value = '''009098'
$q4.Cells.Item(6,3).Value = '景顺长城价值领航两年持有期混合'
$q4.Cells.Item(6,4).Value = '''7.08'
$q4.Cells.Item(6,5).Value = '''92.98'
$q4.Cells.Item(6,6).Value = '''8.03'
$q4.Cells.Item(6,7).Value = '''0.5685'
$q4.Cells.Item(6,8).Value = 4

# row 7: 008060
$refSheet.Cells.Item(7,1).Copy($q4.Cells.Item(7,1))
$q4.Cells.Item(7,1).Value = 5
$q4.Cells.Item(7,2).Value = '''008060'
$q4.Cells.Item(7,3).Value = '景顺长城价值边际灵活配置混合A'
$q4.Cells.Item(7,4).Value = '''6.42'
$q4.Cells.Item(7,5).Value = '''91.73'
$q4.Cells.Item(7,6).Value = '''8.08'
$q4.Cells.Item(7,7).Value = '''0.5187'
$q4.Cells.Item(7,8).Value = 3

# row 8: 217024
$refSheet.Cells.Item(8,1).Copy($q4.Cells.Item(8,1))
$q4.Cells.Item(8,1).Value = 6
$q4.Cells.Item(8,2).Value = '''217024'
$q4.Cells.Item(8,3).Value = '招商安盈债券A'
$q4.Cells.Item(8,4).Value = '''40.95'
$q4.Cells.Item(8,5).Value = '''20.32'
$q4.Cells.Item(8,6).Value = '''0.86'
$q4.Cells.Item(8,7).Value = '''0.3522'
$q4.Cells.Item(8,8).Value = 9

# row 9: 011046
$refSheet.Cells.Item(9,1).Copy($q4.Cells.Item(9,1))
$q4.Cells.Item(9,1).Value = 7
$q4.Cells.Item(9,2).Value = '''011046'
$q4.Cells.Item(9,3).Value = '富国优质企业混合A'
$q4.Cells.Item(9,4).Value = '''6.56'
$q4.Cells.Item(9,5).Value = '''79.84'
$q4.Cells.Item(9,6).Value = '''3.80'
$q4.Cells.Item(9,7).Value = '''0.2493'
$q4.Cells.Item(9,8).Value = 7

# row 10: 014887
$refSheet.Cells.Item(10,1).Copy($q4.Cells.Item(10,1))
$q4.Cells.Item(10,1).Value = 8
$q4.Cells.Item(10,2).Value = '''014887'
$q4.Cells.Item(10,3).Value = '招商安福1年定期开放债券'
$q4.Cells.Item(10,4).Value = '''17.78'
$q4.Cells.Item(10,5).Value = '''33.59'
$q4.Cells.Item(10,6).Value = '''1.17'
$q4.Cells.Item(10,7).Value = '''0.2080'
$q4.Cells.Item(10,8).Value = 9

# row 11: 015779
$refSheet.Cells.Item(11,1).Copy($q4.Cells.Item(11,1))
$q4.Cells.Item(11,1).Value = 9
$q4.Cells.Item(11,2).Value = '''015779'
$q4.Cells.Item(11,3).Value = '景顺长城价值边际灵活配置混合C'
$q4.Cells.Item(11,4).Value = '''2.40'
$q4.Cells.Item(11,5).Value = '''91.73'
$q4.Cells.Item(11,6).Value = '''8.08'
$q4.Cells.Item(11,7).Value = '''0.1939'
$q4.Cells.Item(11,8).Value = 3

# row 12: 016513
$refSheet.Cells.Item(12,1).Copy($q4.Cells.Item(12,1))
$q4.Cells.Item(12,1).Value = 10
$q4.Cells.Item(12,2).Value = '''016513'
$q4.Cells.Item(12,3).Value = '招商安嘉债券'
$q4.Cells.Item(12,4).Value = '''16.27'
$q4.Cells.Item(12,5).Value = '''20.17'
$q4.Cells.Item(12,6).Value = '''0.69'
$q4.Cells.Item(12,7).Value = '''0.1123'
$q4.Cells.Item(12,8).Value = 10

# row 13: 009782
$refSheet.Cells.Item(13,1).Copy($q4.Cells.Item(13,1))
$q4.Cells.Item(13,1).Value = 11
$q4.Cells.Item(13,2).Value = '''009782'
$q4.Cells.Item(13,3).Value = '富国兴泉回报12个月持有期混合A'
$q4.Cells.Item(13,4).Value = '''3.09'
$q4.Cells.Item(13,5).Value = '''74.76'
$q4.Cells.Item(13,6).Value = '''3.43'
$q4.Cells.Item(13,7).Value = '''0.1060'
$q4.Cells.Item(13,8).Value = 5

# row 14: 010029
$refSheet.Cells.Item(14,1).Copy($q4.Cells.Item(14,1))
$q4.Cells.Item(14,1).Value = 12
$q4.Cells.Item(14,2).Value = '''010029'
$q4.Cells.Item(14,3).Value = '富国稳进回报12个月持有期混合A'
$q4.Cells.Item(14,4).Value = '''6.03'
$q4.Cells.Item(14,5).Value = '''20.60'
$q4.Cells.Item(14,6).Value = '''1.63'
$q4.Cells.Item(14,7).Value = '''0.0983'
$q4.Cells.Item(14,8).Value = 3

# row 15: 017090
$refSheet.Cells.Item(15,1).Copy($q4.Cells.Item(15,1))
$q4.Cells.Item(15,1).Value = 13
$q4.Cells.Item(15,2).Value = '''017090'
$q4.Cells.Item(15,3).Value = '景顺长城能源基建混合C'
$q4.Cells.Item(15,4).Value = '''1.00'
$q4.Cells.Item(15,5).Value = '''85.17'
$q4.Cells.Item(15,6).Value = '''8.76'
$q4.Cells.Item(15,7).Value = '''0.0876'
$q4.Cells.Item(15,8).Value = 2

# row 16: 009783
$refSheet.Cells.Item(16,1).Copy($q4.Cells.Item(16,1))
$q4.Cells.Item(16,1).Value = 14
$q4.Cells.Item(16,2).Value = '''009783'
$q4.Cells.Item(16,3).Value = '富国兴泉回报12个月持有期混合C'
$q4.Cells.Item(16,4).Value = '''2.02'
$q4.Cells.Item(16,5).Value = '''74.76'
$q4.Cells.Item(16,6).Value = '''3.43'
$q4.Cells.Item(16,7).Value = '''0.0693'
$q4.Cells.Item(16,8).Value = 5

# row 17: 005732
$refSheet.Cells.Item(17,1).Copy($q4.Cells.Item(17,1))
$q4.Cells.Item(17,1).Value = 15
$q4.Cells.Item(17,2).Value = '''005732'
$q4.Cells.Item(17,3).Value = '富国臻选成长灵活配置混合'
$q4.Cells.Item(17,4).Value = '''1.87'
$q4.Cells.Item(17,5).Value = '''76.07'
$q4.Cells.Item(17,6).Value = '''3.48'
$q4.Cells.Item(17,7).Value = '''0.0651'
$q4.Cells.Item(17,8).Value = 4

# row 18: 009840
$refSheet.Cells.Item(18,1).Copy($q4.Cells.Item(18,1))
$q4.Cells.Item(18,1).Value = 16
$q4.Cells.Item(18,2).Value = '''009840'
$q4.Cells.Item(18,3).Value = '西藏东财量化精选混合A'
$q4.Cells.Item(18,4).Value = '''1.37'
$q4.Cells.Item(18,5).Value = '''84.25'
$q4.Cells.Item(18,6).Value = '''4.40'
$q4.Cells.Item(18,7).Value = '''0.0603'
$q4.Cells.Item(18,8).Value = 4

# row 19: 080005
$refSheet.Cells.Item(19,1).Copy($q4.Cells.Item(19,1))
$q4.Cells.Item(19,1).Value = 17
$q4.Cells.Item(19,2).Value = '''080005'
$q4.Cells.Item(19,3).Value = '长盛量化红利混合'
$q4.Cells.Item(19,4).Value = '''1.66'
$q4.Cells.Item(19,5).Value = '''72.88'
$q4.Cells.Item(19,6).Value = '''2.30'
$q4.Cells.Item(19,7).Value = '''0.0382'
$q4.Cells.Item(19,8).Value = 6

# row 20: 009841
$refSheet.Cells.Item(20,1).Copy($q4.Cells.Item(20,1))
$q4.Cells.Item(20,1).Value = 18
$q4.Cells.Item(20,2).Value = '''009841'
$q4.Cells.Item(20,3).Value = '西藏东财量化精选混合C'
$q4.Cells.Item(20,4).Value = '''0.60'
$q4.Cells.Item(20,5).Value = '''84.25'
$q4.Cells.Item(20,6).Value = '''4.40'
$q4.Cells.Item(20,7).Value = '''0.0264'
$q4.Cells.Item(20,8).Value = 4

# row 21: 011047
$refSheet.Cells.Item(21,1).Copy($q4.Cells.Item(21,1))
$q4.Cells.Item(21,1).Value = 19
$q4.Cells.Item(21,2).Value = '''011047'
$q4.Cells.Item(21,3).Value = '富国优质企业混合C'
$q4.Cells.Item(21,4).Value = '''0.41'
$q4.Cells.Item(21,5).Value = '''79.84'
$q4.Cells.Item(21,6).Value = '''3.80'
$q4.Cells.Item(21,7).Value = '''0.0156'
$q4.Cells.Item(21,8).Value = 7

# row 22: 010030
$refSheet.Cells.Item(22,1).Copy($q4.Cells.Item(22,1))
$q4.Cells.Item(22,1).Value = 20
$q4.Cells.Item(22,2).Value = '''010030'
$q4.Cells.Item(22,3).Value = '富国稳进回报12个月持有期混合C'
$q4.Cells.Item(22,4).Value = '''0.89'
$q4.Cells.Item(22,5).Value = '''20.60'
$q4.Cells.Item(22,6).Value = '''1.63'
$q4.Cells.Item(22,7).Value = '''0.0145'
$q4.Cells.Item(22,8).Value = 3

# row 23: 009514
$refSheet.Cells.Item(23,1).Copy($q4.Cells.Item(23,1))
$q4.Cells.Item(23,1).Value = 21
$q4.Cells.Item(23,2).Value = '''009514'
$q4.Cells.Item(23,3).Value = '创金合信同顺创业板精选股票C'
$q4.Cells.Item(23,4).Value = '''0.15'
$q4.Cells.Item(23,5).Value = '''91.10'
$q4.Cells.Item(23,6).Value = '''1.42'
$q4.Cells.Item(23,7).Value = '''0.0021'
$q4.Cells.Item(23,8).Value = 10

# row 24: 009513
$refSheet.Cells.Item(24,1).Copy($q4.Cells.Item(24,1))
$q4.Cells.Item(24,1).Value = 22
$q4.Cells.Item(24,2).Value = '''009513'
$q4.Cells.Item(24,3).Value = '创金合信同顺创业板精选股票A'
$q4.Cells.Item(24,4).Value = '''0.09'
$q4.Cells.Item(24,5).Value = '''91.10'
$q4.Cells.Item(24,6).Value = '''1.42'
$q4.Cells.Item(24,7).Value = '''0.0013'
$q4.Cells.Item(24,8).Value = 10

# row 25: 012233
$refSheet.Cells.Item(25,1).Copy($q4.Cells.Item(25,1))
$q4.Cells.Item(25,1).Value = 23
$q4.Cells.Item(25,2).Value = '''012233'
$q4.Cells.Item(25,3).Value = '招商安盈债券C'
$q4.Cells.Item(25,4).Value = '''0.01'
$q4.Cells.Item(25,5).Value = '''20.32'
$q4.Cells.Item(25,6).Value = '''0.86'
$q4.Cells.Item(25,7).Value = '''0.0001'
$q4.Cells.Item(25,8).Value = 9

# ---- 2. Update the '总计' (summary) sheet: insert the 2022-Q4 record ----
$zj = $wb.Worksheets.Item('总计')

# Give the new index cell (A10) the same style as the existing index column
# before writing the final table top-to-bottom (row 2 = newest quarter).
$zj.Cells.Item(9,1).Copy($zj.Cells.Item(10,1))

$zj.Cells.Item(2,1).Value = 0
$zj.Cells.Item(2,2).Value = '2022-Q4'
$zj.Cells.Item(2,3).Value = 24
$zj.Cells.Item(2,4).Value = 8.9

$zj.Cells.Item(3,1).Value = 1
$zj.Cells.Item(3,2).Value = '2022-Q3'
$zj.Cells.Item(3,3).Value = 23
$zj.Cells.Item(3,4).Value = 8.19

$zj.Cells.Item(4,1).Value = 2
$zj.Cells.Item(4,2).Value = '2022-Q2'
$zj.Cells.Item(4,3).Value = 26
$zj.Cells.Item(4,4).Value = 6.42

$zj.Cells.Item(5,1).Value = 3
$zj.Cells.Item(5,2).Value = '2022-Q1'
$zj.Cells.Item(5,3).Value = 15
$zj.Cells.Item(5,4).Value = 5.12

$zj.Cells.Item(6,1).Value = 4
$zj.Cells.Item(6,2).Value = '2021-Q4'
$zj.Cells.Item(6,3).Value = 12
$zj.Cells.Item(6,4).Value = 4.16

$zj.Cells.Item(7,1).Value = 5
$zj.Cells.Item(7,2).Value = '2021-Q3'
$zj.Cells.Item(7,3).Value = 61
$zj.Cells.Item(7,4).Value = 23.94

$zj.Cells.Item(8,1).Value = 6
$zj.Cells.Item(8,2).Value = '2021-Q2'
$zj.Cells.Item(8,3).Value = 23
$zj.Cells.Item(8,4).Value = 4.48

$zj.Cells.Item(9,1).Value = 7
$zj.Cells.Item(9,2).Value = '2021-Q1'
$zj.Cells.Item(9,3).Value = 19
$zj.Cells.Item(9,4).Value = 3.05

$zj.Cells.Item(10,1).Value = 8
$zj.Cells.Item(10,2).Value = '2020-Q4'
$zj.Cells.Item(10,3).Value = 11
$zj.Cells.Item(10,4).Value = 1.83

